$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates: force text type to match original inlineStr formatting ---
# (values like "234.57" or "6.04" would otherwise be auto-parsed as numbers by Excel)
$dCells = @("D2","D3","D5","D7","D10","D13","D15","D16","D18","D19","D22","D23","D25","D27","D28","D29","D31","D33","D36","D37","D39","D41","D43","D44","D46","D48")
foreach ($ref in $dCells) { $ws.Range($ref).NumberFormat = "@" }

$ws.Range("D2").Value = '37.526.60'
$ws.Range("D3").Value = '2.075.26'
$ws.Range("D5").Value = '234.57'
$ws.Range("D7").Value = '58.26'
$ws.Range("D10").Value = '59.22'
$ws.Range("D13").Value = '2.381.69'
$ws.Range("D15").Value = '21.04'
$ws.Range("D16").Value = '0.781'
$ws.Range("D18").Value = '2.077.34'
$ws.Range("D19").Value = '37.710.82'
$ws.Range("D22").Value = '0.0₃0816'
$ws.Range("D23").Value = '226.17'
$ws.Range("D25").Value = '2.49'
$ws.Range("D27").Value = '166.36'
$ws.Range("D28").Value = '9.04'
$ws.Range("D29").Value = '1.49'
$ws.Range("D31").Value = '19.26'
$ws.Range("D33").Value = '4.52'
$ws.Range("D36").Value = '4.58'
$ws.Range("D37").Value = '6.04'
$ws.Range("D39").Value = '3.35'
$ws.Range("D41").Value = '4.67'
$ws.Range("D43").Value = '0.0953'
$ws.Range("D44").Value = '1.468.06'
$ws.Range("D46").Value = '95.65'
$ws.Range("D48").Value = '15.84'

# restore default (un-styled) cell style now that the text is committed as a string
foreach ($ref in $dCells) { $ws.Range($ref).Style = "Normal" }

# --- Other column (B, C, E) updates: plain text values, no numeric ambiguity ---
$ws.Range("E2").Value = '  +2.30%  '
$ws.Range("E3").Value = '  +3.20%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("E7").Value = '  +6.32%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +4.18%  '
$ws.Range("E10").Value = '  +3.11%  '
$ws.Range("E11").Value = '  +2.27%  '
$ws.Range("E12").Value = '  -0.46%  '
$ws.Range("E13").Value = '  +3.42%  '
$ws.Range("E14").Value = '  +2.52%  '
$ws.Range("E15").Value = '  +4.55%  '
$ws.Range("E16").Value = '  +3.29%  '
$ws.Range("E17").Value = '  +2.00%  '
$ws.Range("E18").Value = '  +3.14%  '
$ws.Range("E19").Value = '  +3.13%  '
$ws.Range("E20").Value = '  +16.63%  '
$ws.Range("E21").Value = '  +2.77%  '
$ws.Range("E22").Value = '  +2.51%  '
$ws.Range("E23").Value = '  +2.29%  '
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("E25").Value = '  +4.33%  '
$ws.Range("E26").Value = '  +0.55%  '
$ws.Range("E27").Value = '  +2.00%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("E28").Value = '  +4.59%  '
$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("E29").Value = '  +6.74%  '
$ws.Range("E30").Value = '  +2.56%  '
$ws.Range("E31").Value = '  +2.59%  '
$ws.Range("E32").Value = '  +1.66%  '
$ws.Range("E33").Value = '  +3.76%  '
$ws.Range("E34").Value = '  +4.01%  '
$ws.Range("E35").Value = '  +7.03%  '
$ws.Range("E36").Value = '  +8.13%  '
$ws.Range("E37").Value = '  +5.07%  '
$ws.Range("E38").Value = '  +0.05%  '
$ws.Range("E39").Value = '  +1.45%  '
$ws.Range("E40").Value = '  -0.59%  '
$ws.Range("E41").Value = '  +20.78%  '
$ws.Range("E42").Value = '  -1.09%  '
$ws.Range("E43").Value = '  +2.52%  '
$ws.Range("E44").Value = '  +0.72%  '
$ws.Range("E45").Value = '  +6.93%  '
$ws.Range("E46").Value = '  +6.46%  '
$ws.Range("E47").Value = '  +4.13%  '
$ws.Range("E48").Value = '  +3.60%  '
$ws.Range("E49").Value = '  +3.78%  '
$ws.Range("E50").Value = '  +5.54%  '
$ws.Range("E51").Value = '  +1.91%  '
